$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 262, shifting existing rows 262:341 down to 263:342
$ws.Rows.Item(262).Insert()

# Fill the newly inserted row 262 with the new record's data.
# Columns A, B, C, E, F, G, I, O, R carry the same values the row above (old
# row 262, now row 263) has, matching the rest of this "Feria Lagunitas de
# Puerto Montt" / Melon block.
$ws.Range("A262").Value = 4
$ws.Range("B262").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C262").Value = 'Los Lagos'
$ws.Range("D262").Value = "12/22/2022"
$ws.Range("E262").Value = 10
$ws.Range("F262").Value = 100112027
$ws.Range("G262").Value = 'Melón'
$ws.Range("H262").Value = 'Tuna'
$ws.Range("I262").Value = 'Extra'
$ws.Range("J262").Value = 150
$ws.Range("K262").Value = 21000
$ws.Range("L262").Value = 21000
$ws.Range("M262").Value = 21000
$ws.Range("N262").Value = '$/caja 12 unidades'
$ws.Range("O262").Value = "Región de O'Higgins"
$ws.Range("P262").Value = 1750
$ws.Range("Q262").Value = 12
$ws.Range("R262").Value = 'Hortaliza'
